$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.840.06"
$ws.Range("D3").Value = "2.395.71"
$ws.Range("E3").Value = "  +0.87%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "504.94"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.32%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "132.49"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.07%  "
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.551"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.56%  "
$ws.Range("D9").Value = "2.401.98"
$ws.Range("E9").Value = "  +0.20%  "
$ws.Range("E10").Value = "  +0.80%  "
$ws.Range("E11").Value = "  -0.70%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.323"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.12%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.64"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.53%  "
$ws.Range("D14").Value = "2.804.34"
$ws.Range("E14").Value = "  +0.12%  "
$ws.Range("D15").Value = "56.774.24"
$ws.Range("E15").Value = "  +0.67%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.67"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.00%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000133"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.79%  "
$ws.Range("D18").Value = "2.407.82"
$ws.Range("E18").Value = "  +0.23%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.22"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.22%  "
$ws.Range("E20").Value = "  -0.11%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "309.49"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.31%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.27"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.13%  "
$ws.Range("E23").Value = "  -0.09%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.56"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.00%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "67.22"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.06%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.997"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.09%  "
$ws.Range("E27").Value = "  -0.88%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.150"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.91%  "
$ws.Range("E29").Value = "  +2.31%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "176.01"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.80%  "
$ws.Range("E31").Value = "  +1.06%  "
$ws.Range("E32").Value = "  -1.00%  "
$ws.Range("E33").Value = "  +0.88%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.88"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.58%  "
$ws.Range("E35").Value = "  +0.19%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.997"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.19%  "
$ws.Range("E37").Value = "  +0.66%  "
$ws.Range("E38").Value = "  -2.02%  "
$ws.Range("E39").Value = "  +1.27%  "
$ws.Range("B40").Value = "SuiNetwork"
$ws.Range("C40").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.831"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.61%  "
$ws.Range("B41").Value = "OKB"
$ws.Range("C41").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "36.81"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.70%  "
$ws.Range("E42").Value = "  +0.47%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "131.57"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.36%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.38"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.22%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.87"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.59%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.568"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.52%  "
$ws.Range("B47").Value = "Bittensor"
$ws.Range("C47").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "252.14"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.36%  "
$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0911"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.09%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0486"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.57%  "
$ws.Range("E50").Value = "  +1.09%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "17.13"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +7.71%  "
